$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.803.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.030.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.06%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.02%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.023.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.60%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.13%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.126"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.536.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.799.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.032.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.57"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.26"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.47"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.16"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +5.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.31"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.04%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.43"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.74%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.65%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.13%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.56"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0864"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.48%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.16"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +9.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.10"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.40"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.302"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +12.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "396.10"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0359"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.722.17"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.75"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.35"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.06%  "

